$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dataset values (Topsis / No Topsis cost & carbon emissions per year)
$ws.Range("B3").Value = 4115709
$ws.Range("C3").Value = 102554683
$ws.Range("F3").Value = 4507764
$ws.Range("G3").Value = 108503834

$ws.Range("B4").Value = 3610891
$ws.Range("C4").Value = 41626621
$ws.Range("F4").Value = 4358954
$ws.Range("G4").Value = 74426573

$ws.Range("B5").Value = 3257979
$ws.Range("C5").Value = 72636868
$ws.Range("F5").Value = 4284770
$ws.Range("G5").Value = 34509109

$ws.Range("B6").Value = 3541365
$ws.Range("C6").Value = 49609096
$ws.Range("F6").Value = 5389524
$ws.Range("G6").Value = 95089672

$ws.Range("B7").Value = 3874720
$ws.Range("C7").Value = 69924314
$ws.Range("F7").Value = 5101151
$ws.Range("G7").Value = 61715280

$ws.Range("B8").Value = 4352695
$ws.Range("C8").Value = 31844686
$ws.Range("F8").Value = 4960658
$ws.Range("G8").Value = 54521452

$ws.Range("B9").Value = 4254041
$ws.Range("C9").Value = 53078409
$ws.Range("F9").Value = 5289381
$ws.Range("G9").Value = 77853687

$ws.Range("B10").Value = 3867912
$ws.Range("C10").Value = 40875395
$ws.Range("F10").Value = 4857609
$ws.Range("G10").Value = 61853346

$ws.Range("B11").Value = 4717384
$ws.Range("C11").Value = 47643565
$ws.Range("F11").Value = 6295486
$ws.Range("G11").Value = 111963711

$ws.Range("B12").Value = 4997564
$ws.Range("C12").Value = 29411067
$ws.Range("F12").Value = 5180814
$ws.Range("G12").Value = 54556791

$ws.Range("B13").Value = 5715050
$ws.Range("C13").Value = 41315167
$ws.Range("F13").Value = 6746906
$ws.Range("G13").Value = 115129283

$ws.Range("B14").Value = 4389765
$ws.Range("C14").Value = 28956417
$ws.Range("F14").Value = 5682924
$ws.Range("G14").Value = 77300851

$ws.Range("B15").Value = 5385026
$ws.Range("C15").Value = 46170317
$ws.Range("F15").Value = 6093110
$ws.Range("G15").Value = 76329555

$ws.Range("B16").Value = 4989551
$ws.Range("C16").Value = 30004642
$ws.Range("F16").Value = 6938642
$ws.Range("G16").Value = 138404367

$ws.Range("B17").Value = 6464877
$ws.Range("C17").Value = 70230229
$ws.Range("F17").Value = 5248071
$ws.Range("G17").Value = 23414506

$ws.Range("B18").Value = 4108288
$ws.Range("C18").Value = 25511577
$ws.Range("F18").Value = 6870702
$ws.Range("G18").Value = 158040286

# Updated summary text labels
$ws.Range("D24").Value = "Topsis decreased cost by 16%"
$ws.Range("D25").Value = "but TOPSIS got 56% reduction in total CE over the years"

# Recalculate formulas so dependent sums/percentages update
$excel.Calculate()

# Update the active cell selection
$ws.Range("L25").Select()
